# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# --- Rename existing header labels ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsForecast.Move($null, $wsMonthly)

# Re-resolve the sheet handle by name: Move() invalidates the old positional
# reference, so grab it fresh before writing any data into it.
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the header styling used on the other sheets' header row (bold,
# thin box border, centered horizontal / top vertical alignment).
$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$data = @(
    @(44934.99999999999, 6, -0.7427556341377912, 13.08447272778536),
    @(44941.99999999999, 6, -1.530908584557705, 13.80151533585393),
    @(44948.99999999999, 6, -0.5559867935009757, 13.26912925117452),
    @(44976.99999999999, 7, -0.6030283859063142, 13.76208840365125),
    @(44983.99999999999, 7, 0.2039714896898608, 13.58009763564213),
    @(44990.99999999999, 7, -0.3664395164098342, 13.30813350607357),
    @(44997.99999999999, 7, -0.04297871769377699, 14.38810345462439),
    @(45004.99999999999, 7, 0.129650002471256, 14.16448849840646),
    @(45011.99999999999, 7, 0.04144177290754125, 14.51687542271166),
    @(45088.99999999999, 8, 1.696999620240295, 15.66803794764647),
    @(45109.99999999999, 9, 2.324505352382803, 15.73732081493198),
    @(45116.99999999999, 9, 0.8779843536244124, 15.63546259516217),
    @(45123.99999999999, 9, 1.676886761852477, 15.4692662783174),
    @(45130.99999999999, 9, 1.901580622113563, 15.9244615043123),
    @(45137.99999999999, 9, 2.745786039870019, 16.12796648796947),
    @(45151.99999999999, 9, 2.344239507004171, 16.67218479524253),
    @(45207.99999999999, 10, 3.374252754802591, 17.58885419435795),
    @(45221.99999999999, 10, 2.638835053528219, 16.67720254644237),
    @(45242.99999999999, 11, 3.715206566089994, 17.71845040907724),
    @(45277.99999999999, 11, 4.175406189615223, 17.74956126943313),
    @(45298.99999999999, 11, 4.429521633836145, 17.87731632874363),
    @(45312.99999999999, 12, 4.458937000880901, 19.12302946918552),
    @(45319.99999999999, 12, 5.130542228514225, 19.19773057665909),
    @(45326.99999999999, 12, 4.59095060839462, 19.21958799911343),
    @(45333.99999999999, 12, 4.490600740028049, 19.25076626778501),
    @(45340.99999999999, 12, 4.383465055690807, 18.8657176385625),
    @(45347.99999999999, 12, 5.075719795434793, 19.50029963200099),
    @(45354.99999999999, 12, 5.71157472639438, 18.82929784701776),
    @(45361.99999999999, 12, 5.871674868699605, 19.67160482928913),
    @(45368.99999999999, 12, 5.619700204520007, 19.2262546771913)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Column A (the "ds" date column) uses the same custom date format as the
# date columns on the other two sheets.
$wsForecast.Range("A2:A31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
